$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 3.3
$ws.Range("I2").Value = 2.4
$ws.Range("AN2").Value = 5
$ws.Range("AO2").Value = 19
$ws.Range("AQ2").Value = 67
$ws.Range("AR2").Value = 101

# Row 6 updates
$ws.Range("J6").Value = 3.2
$ws.Range("N6").Value = 8
$ws.Range("AA6").Value = 21
$ws.Range("AH6").Value = 8
$ws.Range("AK6").Value = 34
$ws.Range("AN6").Value = 4.33
